$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The ratio formula in D4 is no longer needed - remove it
$ws.Range("D4").ClearContents()

# Add new measurement rows
$ws.Range("A5").Value = "all black + black circle"
$ws.Range("C5").Value = 20614.9

$ws.Range("A6").Value = "all black + black circle + grease a lot"
$ws.Range("C6").Value = 22001.7

# Row 7 is intentionally left blank

$ws.Range("A8").Value = "all black 2 + grease a lot"
$ws.Range("C8").Value = 22321.5

$ws.Range("A9").Value = "all black 2 + grease a lot + not center"

# Widen column A to fit the new, longer labels
$ws.Columns.Item(1).ColumnWidth = 32.65

# Leave the selection where it was left after entering the data
$ws.Range("E22").Select()
